# Update workbook metadata: absPath (best-effort; Excel COM doesn't expose this
# directly, so we skip it - it is not something reachable through the object
# model in a meaningful way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows (2-5) with new values ---
$ws.Range("B2").Value = 49434
$ws.Range("C2").Value = 49741

$ws.Range("B3").Value = 49546
$ws.Range("C3").Value = 49741

$ws.Range("B4").Value = 49546
$ws.Range("C4").Value = 49870

$ws.Range("B5").Value = 49434
$ws.Range("C5").Value = 49870

# Row 6 stays the same (37.5, 51499.378, 49805.856, 0)

# --- Add new block starting at row 7: header row (shared strings) ---
$ws.Range("A7").Value = "반지름"
$ws.Range("B7").Value = "중심점X"
$ws.Range("C7").Value = "중심점Y"
$ws.Range("D7").Value = "중심점Z"

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 39434
$ws.Range("C8").Value = 39620
$ws.Range("D8").Value = 0

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 39546
$ws.Range("C9").Value = 39620
$ws.Range("D9").Value = 0

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 39546
$ws.Range("C10").Value = 39870
$ws.Range("D10").Value = 0

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 39434
$ws.Range("C11").Value = 39870
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = 37.5
$ws.Range("B12").Value = 51499.377999999997
$ws.Range("C12").Value = 49805.856
$ws.Range("D12").Value = 0

# --- Add third block starting at row 13: header row (shared strings) ---
$ws.Range("A13").Value = "반지름"
$ws.Range("B13").Value = "중심점X"
$ws.Range("C13").Value = "중심점Y"
$ws.Range("D13").Value = "중심점Z"

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = 29434
$ws.Range("C14").Value = 29661
$ws.Range("D14").Value = 0

$ws.Range("A15").Value = 8
$ws.Range("B15").Value = 29546
$ws.Range("C15").Value = 29661
$ws.Range("D15").Value = 0

$ws.Range("A16").Value = 8
$ws.Range("B16").Value = 29546
$ws.Range("C16").Value = 29870
$ws.Range("D16").Value = 0

$ws.Range("A17").Value = 8
$ws.Range("B17").Value = 29434
$ws.Range("C17").Value = 29870
$ws.Range("D17").Value = 0

$ws.Range("A18").Value = 37.5
$ws.Range("B18").Value = 51499.377999999997
$ws.Range("C18").Value = 49805.856
$ws.Range("D18").Value = 0

# --- Update the selection / active cell to match the edited state ---
$ws.Range("F15").Select()
